# semana 31 de 2024
# Update Esperado (C), Observado (D) and valor p (E) figures on the
# poisson data sheet for the week-31/2024 refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 0.27

$ws.Range("C4").Value = 1
$ws.Range("E4").Value = 0.37

$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 7
$ws.Range("E5").Value = 0.1

$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 90

$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = 0.18

$ws.Range("C11").Value = 40
$ws.Range("D11").Value = 30
$ws.Range("E11").Value = 0.02

$ws.Range("D13").Value = 4
$ws.Range("E13").Value = 0.02

$ws.Range("C14").Value = 5
$ws.Range("D14").Value = 5
$ws.Range("E14").Value = 0.18

$ws.Range("C15").Value = 15
$ws.Range("D15").Value = 0

$ws.Range("C17").Value = 2
$ws.Range("E17").Value = 0.14

$ws.Range("C18").Value = 0
$ws.Range("E18").Value = 0

$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 0.1

$ws.Range("C21").Value = 7
$ws.Range("D21").Value = 3
$ws.Range("E21").Value = 0.05

$ws.Range("D24").Value = 3
$ws.Range("E24").Value = 0.06

$ws.Range("C27").Value = 5
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = 0.08

$ws.Range("C28").Value = 2
$ws.Range("E28").Value = 0.27

$ws.Range("D29").Value = 1
$ws.Range("E29").Value = 0

$ws.Range("C30").Value = 3
$ws.Range("E30").Value = 0.05

$ws.Range("C33").Value = 8
$ws.Range("D33").Value = 4
$ws.Range("E33").Value = 0.06

$ws.Range("C34").Value = 9
$ws.Range("D34").Value = 1
$ws.Range("E34").Value = 0

$ws.Range("C35").Value = 8
$ws.Range("E35").Value = 0.14
